# Commit message: "removed domain part from sheet name. Made generic for
# all domains. Any input file should have the below sheet name input
# files: Requirements Domain Dataelements"
#
# The workbook shipped with a domain-specific sheet name
# ("Banking-Requirements"); rename it to the generic "Requirements" so the
# downstream script can rely on a fixed, domain-agnostic sheet name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Requirements"
